$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2..17 across columns D, L, M, N, O, P, S.
# (The underlying edit re-orders the weekly price records; values below
# reflect the final state shown in the target diff.)
$rows = @{
    2  = @{ D = 44446; L = "Primera";  M = 60;  N = 21000; O = 22000; P = 21500; S = 2150 }
    3  = @{ D = 44848; L = "Especial"; M = 60;  N = 24000; O = 25000; P = 24500; S = 2450 }
    4  = @{ D = 44848; L = "Primera";  M = 120; N = 21000; O = 22000; P = 21500; S = 2150 }
    5  = @{ D = 44868; L = "Especial"; M = 60;  N = 26000; O = 26000; P = 26000; S = 2600 }
    6  = @{ D = 44839; L = "Primera";  M = 120; N = 25000; O = 26000; P = 25500; S = 2550 }
    7  = @{ D = 44874; L = "Especial"; M = 30;  N = 25000; O = 25000; P = 25000; S = 2500 }
    8  = @{ D = 44874; L = "Primera";  M = 80;  N = 23000; O = 24000; P = 23500; S = 2350 }
    9  = @{ D = 44461; L = "Especial"; M = 60;  N = 31000; O = 32000; P = 31500; S = 3150 }
    10 = @{ D = 44461; L = "Primera";  M = 30;  N = 30000; O = 30000; P = 30000; S = 3000 }
    11 = @{ D = 44841; L = "Primera";  M = 60;  N = 23000; O = 24000; P = 23500; S = 2350 }
    12 = @{ D = 44448; L = "Primera";  M = 60;  N = 21000; O = 22000; P = 21500; S = 2150 }
    13 = @{ D = 44487; L = "Primera";  M = 30;  N = 23000; O = 24000; P = 23500; S = 2350 }
    14 = @{ D = 44452; L = "Primera";  M = 60;  N = 21000; O = 22000; P = 21500; S = 2150 }
    15 = @{ D = 44460; L = "Especial"; M = 60;  N = 31000; O = 32000; P = 31500; S = 3150 }
    16 = @{ D = 44460; L = "Primera";  M = 30;  N = 30000; O = 30000; P = 30000; S = 3000 }
    17 = @{ D = 44447; L = "Primera";  M = 60;  N = 21000; O = 22000; P = 21500; S = 2150 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($r, 12).Value = $vals.L   # L - Calidad
    $ws.Cells.Item($r, 13).Value = $vals.M   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $vals.N   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals.O   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals.P   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $vals.S   # S - Precio $/Kg
}
